$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Investment_billion_USD"
$ws.Range("B2").Select()
